$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the trailing "_GoBack" bookmark in the last paragraph and
#    the (now orphaned) run that used to sit right after it holding a
#    single space. Both the bookmark and that run get relocated into
#    the brand-new paragraph that is added further below (it ends up
#    becoming the "e" in "save").
# ------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

$bookmark = $d.Bookmarks("_GoBack")
$oldRunStart = $bookmark.Start
$oldRunEnd = $lastPara.Range.End - 1   # stop right before the paragraph mark

# Delete the bookmark itself first (collapses to nothing, does not
# touch surrounding text).
$bookmark.Delete()

# Now remove the leftover single-space run that used to trail the
# bookmark - its content is being replaced by new text further down.
if ($oldRunEnd -gt $oldRunStart) {
    $oldRun = $d.Range($oldRunStart, $oldRunEnd)
    $oldRun.Delete()
}

# ------------------------------------------------------------------
# 2) Append a brand-new (plain, non-bold) run holding a single space
#    right after the bold "&#8364;" run, finishing the paragraph.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfParagraph = $lastPara.Range.End
$tailRange = $d.Range($endOfParagraph, $endOfParagraph)
$tailRange.InsertAfter(" ")

# ------------------------------------------------------------------
# 3) Insert the two new paragraphs after it:
#      - an empty "Prrafodelista" paragraph (ind left=1080, no bullet)
#      - a numbered "Prrafodelista" list item with the new note text,
#        a bold "no", and the relocated "_GoBack" bookmark sitting
#        between "sav" and "e" (i.e. in the middle of the word "save").
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $lastPara.Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$newParagraphsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:ind w:left="1080"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Los momentos que haya que inicializar con el momento de creaci&#243;n de la entidad, deben inicializarse en el m&#233;todo </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>create</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> y </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
              </w:rPr>
              <w:t>no</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> en el </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>sav</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:t>e</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($newParagraphsXml) | Out-Null
